# Compare admittance spectrum of GFM and GFL
# Applies the authoring changes to SgInfiniteBus.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Bus": tweak bus parameters, move selection
# ---------------------------------------------------------------------
$bus = $wb.Worksheets.Item("Bus")
$bus.Range("C4").Value = 1
$bus.Range("I4").Value = -999
$bus.Range("J4").Value = 999
$bus.Range("E5").Value = 0.5
$bus.Range("F5").Value = 0
$bus.Range("I5").Value = -999
$bus.Range("J5").Value = 999
$bus.Range("I6").Select()

# ---------------------------------------------------------------------
# Sheet "Apparatus": insert a header row for the GFM columns and
# populate the second apparatus (bus 2) with the new parameters.
# ---------------------------------------------------------------------
$app = $wb.Worksheets.Item("Apparatus")
$app.Rows.Item(3).Insert()

$app.Range("C3").Value = "J (pu)"
$app.Range("D3").Value = "D (pu)"
$app.Range("E3").Value = "wL (pu)"
$app.Range("F3").Value = "R (pu)"

$app.Range("B6").Value = 1
$app.Range("C6").Value = 1
$app.Range("D6").Value = 5
$app.Range("E6").Value = 0.05
$app.Range("F6").Formula = "=E6/10"

# ---------------------------------------------------------------------
# Sheet "NetworkLine": derive R from X via formula, update X, select row 5
# ---------------------------------------------------------------------
$netline = $wb.Worksheets.Item("NetworkLine")
$netline.Range("D4").Value = 0.3
$netline.Range("C4").Formula = "=D4/10"
$netline.Rows.Item(5).Select()

# ---------------------------------------------------------------------
# Sheet "Advance": disable the Simulink model flag, move selection
# ---------------------------------------------------------------------
$adv = $wb.Worksheets.Item("Advance")
$adv.Range("B8").Value = 0
$adv.Range("B9").Select()

# ---------------------------------------------------------------------
# Make "Apparatus" the active sheet/tab (must be the last Select call)
# ---------------------------------------------------------------------
$app.Activate()
$app.Range("B6").Select()
